$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213139176368713
$ws.Range("B1").Value = 2.633672952651978
$ws.Range("D1").Value = 2.159224033355713
$ws.Range("E1").Value = 1.15968656539917
